$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the "Examples" column values (new shared strings, in row order so
#    they land at sharedStrings indices 23..32 like the target workbook).
# ---------------------------------------------------------------------------
$ws.Range("C2").Value  = "They announced their engagement at the family dinner."
$ws.Range("C3").Value  = "The marriage ceremony was held in a small, ancient church."
$ws.Range("C4").Value  = "The bride looked beautiful in her white gown."
$ws.Range("C5").Value  = "The groom waited nervously at the altar for the bride."
$ws.Range("C6").Value  = "The graduation ceremony starts at 9:00 AM sharp."
$ws.Range("C7").Value  = "After the wedding, guests went to the hotel for the reception."
$ws.Range("C8").Value  = "The newlyweds are moving into their new apartment next week."
$ws.Range("C9").Value  = "They went to Hawaii for their honeymoon."
$ws.Range("C10").Value = "She introduced her fiancé to her colleagues."
$ws.Range("C11").Value = "He bought a diamond ring for his fiancée."

# ---------------------------------------------------------------------------
# 2. Fonts.
#    Column A (English vocab) + header row -> bold Arial 11 FF1F1F1F
#    Column B/C (explanations + examples)  -> regular Arial 11 FF1F1F1F
#    (Note: multi-area unions only reliably apply to the first area in this
#    engine, so the header row and column A are styled as separate calls.)
# ---------------------------------------------------------------------------
foreach ($addr in @("A1:C1", "A2:A11")) {
    $r = $ws.Range($addr)
    $r.Font.Name  = "Arial"
    $r.Font.Size  = 11
    $r.Font.Bold  = $true
    $r.Font.Color = 2039583   # RGB(0x1F,0x1F,0x1F)
}

$plainRange = $ws.Range("B2:C11")
$plainRange.Font.Name  = "Arial"
$plainRange.Font.Size  = 11
$plainRange.Font.Bold  = $false
$plainRange.Font.Color = 2039583  # RGB(0x1F,0x1F,0x1F)

# ---------------------------------------------------------------------------
# 3. Borders: medium black box around every cell of the table.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:C11")
$tableRange.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$tableRange.Borders.Weight    = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlMedium
$tableRange.Borders.Color     = 0

# ---------------------------------------------------------------------------
# 4. Alignment: left / center-vertical / wrap / indent 1 / left-to-right.
# ---------------------------------------------------------------------------
$tableRange.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$tableRange.VerticalAlignment   = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$tableRange.WrapText            = $true
$tableRange.IndentLevel         = 1
$tableRange.ReadingOrder        = [Microsoft.Office.Interop.Excel.XlReadingOrder]::xlLTR

# ---------------------------------------------------------------------------
# 5. Column widths (closest values this engine's character-width grid can
#    represent to the target bestFit widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.7
$ws.Columns.Item(2).ColumnWidth = 41.7
$ws.Columns.Item(3).ColumnWidth = 61.7

# ---------------------------------------------------------------------------
# 6. Row heights.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight  = 30.75
$ws.Rows.Item(2).RowHeight  = 15.75
$ws.Rows.Item(3).RowHeight  = 30.75
$ws.Rows.Item(4).RowHeight  = 15.75
$ws.Rows.Item(5).RowHeight  = 15.75
$ws.Rows.Item(6).RowHeight  = 15.75
$ws.Rows.Item(7).RowHeight  = 15.75
$ws.Rows.Item(8).RowHeight  = 15.75
$ws.Rows.Item(9).RowHeight  = 15.75
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 7. Selection.
# ---------------------------------------------------------------------------
$ws.Range("E9").Select()
